$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3 so existing A3 ("GDP vs GDI.xlsx") shifts down to A4,
# leaving a blank row 3 to hold the original Bitcoin title text.
$ws.Rows.Item(3).Insert()

# Update A2 with the new title (add "(with extrapolation)")
$ws.Range("A2").Value = "Bitcoin price change is fueled by global monetary growth (with extrapolation).xlsx"

# A3 now holds the original Bitcoin title text
$ws.Range("A3").Value = "Bitcoin price change is fueled by global monetary growth.xlsx"

# A4 already has "GDP vs GDI.xlsx" (shifted down from old A3)

# Add new row 5 with the new title
$ws.Range("A5").Value = "Other deposit liabilities (ODL) shows where US M2 is heading.xlsx"
